$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows for person (household) type 4 with the same wfh_share pattern
# used for the other household types (0, 0.5, 1).
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 0

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 0.5

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 1

$ws.Range("A13").Select()
